$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0 -Millisecond 0).AddDays(45182)

for ($r = 2; $r -le 250; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDate
}
